$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Locoto at Vega Modelo de Temuco.
# It belongs chronologically right after the existing row 31 (date 44795),
# so insert a new row at position 32, which pushes the former rows 32-41
# down to rows 33-42 while keeping their data intact.
$ws.Rows("32").Insert()

# Populate the newly inserted row 32 with the new record's data.
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 44809
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = 100112042
$ws.Range("G32").Value = "Locoto"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 150
$ws.Range("K32").Value = 2700
$ws.Range("L32").Value = 2700
$ws.Range("M32").Value = 2700
$ws.Range("N32").Value = "`$/kilo"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 2700
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"
